$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the "productType" id values in column A: they had been left as 6..10,
# renumber them sequentially starting at 1 (update product logic fix).
for ($i = 0; $i -lt 5; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $i + 1
}

# The "price" column (C) had accidentally picked up an extra/stray number
# format style; clear it back to the workbook's default (Normal) style so
# the header + values go back to being unstyled.
$ws.Columns.Item(3).ClearFormats()

# Move the active selection to C11.
$ws.Range("C11").Select()
